# Adds three new data rows (28, 29, 30) to the NIFTY_Options_Analysis sheet,
# matching the 2026-01-30 entries described by the commit diff.
#
# Column layout (row 1 headers):
#  A Date             B Time              C Signal            D Signal_Tier
#  E Position_Size     F Premium_Quality   G Total_Score       H NIFTY_Spot
#  I VIX               J VIX_Trend         K VIX_Score         L IV_Rank
#  M Market_Regime     N Regime_Score      O OI_Pattern        P OI_Score
#  Q Theta_Score       R Gamma_Score       S Vega_Score        T Best_Strategy
#  U Expiry_1          V Days_To_Expiry_1  W Straddle_Premium  X Straddle_Theta
#  Y Straddle_Gamma    Z Strangle_Premium  AA Strangle_Theta   AB Strangle_Gamma
#  AC Recommendation   AD Risk_Factors     AE Telegram_Sent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that already carry a "0.0" / "0.00" / "0.000000" number format in the
# existing sheet (style indices 5 / 6 / 7 respectively). Column V is numeric
# but intentionally keeps the plain "General" style (same as the text cells).
$fmt1Cols = @("G", "K", "L", "N", "P", "Q", "R", "S")          # 0.0
$fmt2Cols = @("H", "I", "J", "W", "X", "Z", "AA")              # 0.00
$fmt3Cols = @("Y", "AB")                                       # 0.000000
$plainTextCols = @("A", "B", "E", "F", "M", "O", "T", "U", "AC", "AD", "AE")

# A scratch cell well outside the used range. We stage "risky" text there
# (values that Excel's smart-entry would otherwise reinterpret as a date /
# time / percentage / number) under a forced text format, then copy *only
# the value* onto the destination cell. That keeps the destination's own
# number format at "General" (so it lands on the same shared style as its
# neighbours) while still storing a literal text string.
$scratch = $ws.Range("BZ5000")

function Set-PlainCellFormat {
    param($cell)
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.WrapText = $true
    $cell.Borders.LineStyle = 1
}

function Set-SafeText {
    # Writes $text into $cell without letting Excel auto-convert it into a
    # date/time/percentage/number, and without leaving $cell's own
    # NumberFormat changed away from "General".
    param($cell, [string]$text)

    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear()
}

function Set-TextCell {
    # Plain word/sentence values (AVOID, SELL, STRADDLE, free-text notes, ...)
    # never get reinterpreted by Excel, so they can be assigned directly.
    param($cell, [string]$text)
    $cell.Value = $text
}

function Fill-Row {
    param($row, $data, $signal)

    foreach ($col in $data.Keys) {
        $cell = $ws.Range($col + $row)
        Set-PlainCellFormat $cell
    }

    foreach ($col in $fmt1Cols) { $ws.Range($col + $row).NumberFormat = "0.0" }
    foreach ($col in $fmt2Cols) { $ws.Range($col + $row).NumberFormat = "0.00" }
    foreach ($col in $fmt3Cols) { $ws.Range($col + $row).NumberFormat = "0.000000" }

    foreach ($col in $data.Keys) {
        $cell = $ws.Range($col + $row)
        $val = $data[$col]

        if ($col -eq "C" -or $col -eq "D") {
            continue # handled separately below (Signal / Signal_Tier colouring)
        }

        if ($val -is [string]) {
            if (($plainTextCols -contains $col) -and ($col -eq "A" -or $col -eq "B" -or $col -eq "E" -or $col -eq "U")) {
                if ($val -eq "") {
                    # Leave blank (matches an empty inline string cell).
                } else {
                    Set-SafeText $cell $val
                }
            } else {
                Set-TextCell $cell $val
            }
        } else {
            $cell.Value = $val
        }
    }

    # Signal (C) / Signal_Tier (D) colour coding.
    $cSig = $ws.Range("C" + $row)
    $dTier = $ws.Range("D" + $row)
    Set-PlainCellFormat $cSig
    Set-PlainCellFormat $dTier
    Set-TextCell $cSig $data["C"]
    Set-TextCell $dTier $data["D"]

    if ($signal -eq "AVOID") {
        $cSig.Interior.Color = 13551615   # FFC7CE (BGR)
        $cSig.Font.Color = 393372         # 9C0006 (BGR)
        $cSig.Font.Bold = $true

        $dTier.Interior.Color = 255       # FF0000 (BGR)
        $dTier.Font.Color = 16777215      # FFFFFF (BGR)
        $dTier.Font.Bold = $true
    } else {
        # SELL / SELL_STRONG
        $cSig.Interior.Color = 13561798   # C6EFCE (BGR)
        $cSig.Font.Color = 24832          # 006100 (BGR)
        $cSig.Font.Bold = $true

        $dTier.Interior.Color = 5287936   # 00B050 (BGR)
        $dTier.Font.Color = 16777215      # FFFFFF (BGR)
        $dTier.Font.Bold = $true
    }
}

# ---------------------------------------------------------------------
# Row 28 - AVOID (hard veto, trending day)
# ---------------------------------------------------------------------
$row28 = [ordered]@{
    "A" = "2026-01-30"; "B" = "10:00:11"; "C" = "AVOID"; "D" = "AVOID";
    "E" = "100%"; "F" = "TRADEABLE"; "G" = 0; "H" = 25284.7; "I" = 13.84;
    "J" = -0.61; "K" = 0; "L" = 65.59999999999999; "M" = "UNKNOWN"; "N" = 0;
    "O" = "UNKNOWN"; "P" = 0; "Q" = 0; "R" = 0; "S" = 0; "T" = "NONE";
    "U" = ""; "V" = 0; "W" = 0; "X" = 0; "Y" = 0; "Z" = 0; "AA" = 0; "AB" = 0;
    "AC" = "HARD VETO: CPR TRENDING DAY: Price 25284.70 below BC 25308.97 - BEARISH TRENDING DAY likely";
    "AD" = "CPR TRENDING DAY: Price 25284.70 below BC 25308.97 - BEARISH TRENDING DAY likely";
    "AE" = "Yes";
}
Fill-Row 28 $row28 "AVOID"

# ---------------------------------------------------------------------
# Row 29 - SELL / SELL_STRONG (excellent conditions)
# ---------------------------------------------------------------------
$row29 = [ordered]@{
    "A" = "2026-01-30"; "B" = "10:15:12"; "C" = "SELL"; "D" = "SELL_STRONG";
    "E" = "100%"; "F" = "EXCELLENT"; "G" = 80.2; "H" = 25317; "I" = 13.73;
    "J" = -0.72; "K" = 75; "L" = 62.6; "M" = "NEUTRAL"; "N" = 100;
    "O" = "LONG_UNWINDING"; "P" = 70; "Q" = 63.9; "R" = 91; "S" = 90;
    "T" = "STRADDLE"; "U" = "2026-02-10"; "V" = 11; "W" = 537.15; "X" = 31.97;
    "Y" = 0.000904; "Z" = 444.8; "AA" = 31.78; "AB" = 0.000898;
    "AC" = "Excellent conditions for option selling";
    "AD" = "No significant risks identified";
    "AE" = "Yes";
}
Fill-Row 29 $row29 "SELL"

# ---------------------------------------------------------------------
# Row 30 - SELL / SELL_STRONG (excellent conditions)
# ---------------------------------------------------------------------
$row30 = [ordered]@{
    "A" = "2026-01-30"; "B" = "14:00:11"; "C" = "SELL"; "D" = "SELL_STRONG";
    "E" = "100%"; "F" = "EXCELLENT"; "G" = 80.2; "H" = 25344.6; "I" = 13.78;
    "J" = -0.67; "K" = 75; "L" = 64.09999999999999; "M" = "NEUTRAL"; "N" = 100;
    "O" = "SHORT_COVERING"; "P" = 70; "Q" = 63.9; "R" = 91; "S" = 90;
    "T" = "STRADDLE"; "U" = "2026-02-10"; "V" = 11; "W" = 544.45; "X" = 31.97;
    "Y" = 0.000904; "Z" = 452.7; "AA" = 31.79; "AB" = 0.000899;
    "AC" = "Excellent conditions for option selling";
    "AD" = "No significant risks identified";
    "AE" = "Yes";
}
Fill-Row 30 $row30 "SELL"
